$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 402, shifting existing rows 402-415 down to 403-416.
$ws.Rows("402").Insert()

# Populate the newly inserted row 402 with a fresh weekly record (copy of the
# row that used to be at 402, but with an updated date and origin).
$ws.Cells.Item(402, 1).Value = 4
$ws.Cells.Item(402, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(402, 3).Value = "Los Lagos"
$ws.Cells.Item(402, 4).Value = 44939
$ws.Cells.Item(402, 5).Value = 10
$ws.Cells.Item(402, 6).Value = 100114014
$ws.Cells.Item(402, 7).Value = "Betarraga"
$ws.Cells.Item(402, 8).Value = "Sin especificar"
$ws.Cells.Item(402, 9).Value = "Primera"
$ws.Cells.Item(402, 10).Value = 1200
$ws.Cells.Item(402, 11).Value = 1000
$ws.Cells.Item(402, 12).Value = 1000
$ws.Cells.Item(402, 13).Value = 1000
$ws.Cells.Item(402, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(402, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(402, 16).Value = 200
$ws.Cells.Item(402, 17).Value = 5
$ws.Cells.Item(402, 18).Value = "Hortaliza"
